$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1422953928028505
$ws.Range("D2").Value = 0.02986979417712377
$ws.Range("E2").Value = 0.1132051331550912
$ws.Range("F2").Value = 3.863606298255917
$ws.Range("G2").Value = 0.002591977397823942
$ws.Range("J2").Value = 0.2359202244399086
$ws.Range("K2").Value = 2.759205320682497
$ws.Range("N2").Value = 2.601452937013775
$ws.Range("B3").Value = 0.1330337524147467
$ws.Range("D3").Value = 0.02939990356929556
$ws.Range("E3").Value = 0.1110817402023621
$ws.Range("F3").Value = 3.826381375827339
$ws.Range("G3").Value = 0.002598035515932952
$ws.Range("J3").Value = 0.2312728858244952
$ws.Range("K3").Value = 2.628518133706507
$ws.Range("N3").Value = 2.616286902221631
$ws.Range("B4").Value = 0.1274210055145915
$ws.Range("D4").Value = 0.02912535459915588
$ws.Range("E4").Value = 0.1098432858726319
$ws.Range("F4").Value = 3.80562069958961
$ws.Range("G4").Value = 0.002601948534721726
$ws.Range("J4").Value = 0.2285620034034324
$ws.Range("K4").Value = 2.550079440510842
$ws.Range("N4").Value = 2.626138525903201
$ws.Range("B5").Value = 0.1251524408119877
$ws.Range("D5").Value = 0.02901700494232884
$ws.Range("E5").Value = 0.1093549835966563
$ws.Range("F5").Value = 3.797685331502123
$ws.Range("G5").Value = 0.002603591907034501
$ws.Range("J5").Value = 0.2274929756088255
$ws.Range("K5").Value = 2.518565374956609
$ws.Range("N5").Value = 2.630339640690444
$ws.Range("B6").Value = 0.1247768777807465
$ws.Range("D6").Value = 0.02899922742517447
$ws.Range("E6").Value = 0.1092748892266684
$ws.Range("F6").Value = 3.796399302003209
$ws.Range("G6").Value = 0.002603867739347924
$ws.Range("J6").Value = 0.2273176140813007
$ws.Range("K6").Value = 2.513359594506824
$ws.Range("N6").Value = 2.631048483171121
$ws.Range("B7").Value = 0.1273903351307126
$ws.Range("D7").Value = 0.0291238790336763
$ws.Range("E7").Value = 0.1098366342035568
$ws.Range("F7").Value = 3.805511558546613
$ws.Range("G7").Value = 0.002601970499994229
$ws.Range("J7").Value = 0.2285474418918341
$ws.Range("K7").Value = 2.549652611651709
$ws.Range("N7").Value = 2.626194429126357
$ws.Range("B8").Value = 0.1390866797073471
$ws.Range("D8").Value = 0.02970488864748333
$ws.Range("E8").Value = 0.1124594011609936
$ws.Range("F8").Value = 3.850334869225293
$ws.Range("G8").Value = 0.002594026222030558
$ws.Range("J8").Value = 0.2342881139666417
$ws.Range("K8").Value = 2.713767996766023
$ws.Range("N8").Value = 2.606413093235162
$ws.Range("B9").Value = 0.1626075595674905
$ws.Range("D9").Value = 0.03095431094836698
$ws.Range("E9").Value = 0.1181234809871086
$ws.Range("F9").Value = 3.954969090947515
$ws.Range("G9").Value = 0.002579973283323812
$ws.Range("J9").Value = 0.2466861241871641
$ws.Range("K9").Value = 3.050077680984657
$ws.Range("N9").Value = 2.57354228938506
$ws.Range("B10").Value = 0.1802435652080447
$ws.Range("D10").Value = 0.03193851352101973
$ws.Range("E10").Value = 0.1226066486705335
$ws.Range("F10").Value = 4.042211549178404
$ws.Range("G10").Value = 0.002570567432553006
$ws.Range("J10").Value = 0.256504263175259
$ws.Range("K10").Value = 3.306255953653988
$ws.Range("N10").Value = 2.553030338884042
$ws.Range("B11").Value = 0.1883436943725343
$ws.Range("D11").Value = 0.03240047595943452
$ws.Range("E11").Value = 0.1247170105541073
$ws.Range("F11").Value = 4.08418887451694
$ws.Range("G11").Value = 0.002566485565837725
$ws.Range("J11").Value = 0.2611279752003099
$ws.Range("K11").Value = 3.424832262645566
$ws.Range("N11").Value = 2.54449534584387
$ws.Range("B12").Value = 0.1914220860035414
$ws.Range("D12").Value = 0.03257744157168219
$ws.Range("E12").Value = 0.125526421816005
$ws.Range("F12").Value = 4.100416776577816
$ws.Range("G12").Value = 0.002564967997813891
$ws.Range("J12").Value = 0.2629017246886036
$ws.Range("K12").Value = 3.470031693118699
$ws.Range("N12").Value = 2.541378385350214
$ws.Range("B13").Value = 0.1907586092367382
$ws.Range("D13").Value = 0.03253923883305276
$ws.Range("E13").Value = 0.1253516431767494
$ws.Range("F13").Value = 4.096906998801501
$ws.Range("G13").Value = 0.002565293584142288
$ws.Range("J13").Value = 0.2625186962254418
$ws.Range("K13").Value = 3.460283908314295
$ws.Range("N13").Value = 2.542044551853209
$ws.Range("B14").Value = 0.1885967345566115
$ws.Range("D14").Value = 0.03241499443348062
$ws.Range("E14").Value = 0.124783395263151
$ws.Range("F14").Value = 4.085517284594118
$ws.Range("G14").Value = 0.002566360151651536
$ws.Range("J14").Value = 0.2612734434268731
$ws.Range("K14").Value = 3.428544869041048
$ws.Range("N14").Value = 2.544236602302618
$ws.Range("B15").Value = 0.1872739602403612
$ws.Range("D15").Value = 0.03233915505748541
$ws.Range("E15").Value = 0.1244366653994788
$ws.Range("F15").Value = 4.078584071786963
$ws.Range("G15").Value = 0.002567017115269283
$ws.Range("J15").Value = 0.2605136724190089
$ws.Range("K15").Value = 3.409142605992088
$ws.Range("N15").Value = 2.54559429796231
$ws.Range("B16").Value = 0.1797157533811315
$ws.Range("D16").Value = 0.03190860828046738
$ws.Range("E16").Value = 0.12247016435402
$ws.Range("F16").Value = 4.039514544840813
$ws.Range("G16").Value = 0.002570838139820926
$ws.Range("J16").Value = 0.2562052783175375
$ws.Range("K16").Value = 3.298548082643947
$ws.Range("N16").Value = 2.55360418868058
$ws.Range("B17").Value = 0.1750988072958819
$ws.Range("D17").Value = 0.03164811716166582
$ws.Range("E17").Value = 0.1212819935530192
$ws.Range("F17").Value = 4.016135129261244
$ws.Range("G17").Value = 0.002573232525991465
$ws.Range("J17").Value = 0.2536026830396452
$ws.Range("K17").Value = 3.231227025092039
$ws.Range("N17").Value = 2.558722289571051
$ws.Range("B18").Value = 0.172450554316967
$ws.Range("D18").Value = 0.0314996312839213
$ws.Range("E18").Value = 0.1206052631847676
$ws.Range("F18").Value = 4.002903318873365
$ws.Range("G18").Value = 0.002574628257325507
$ws.Range("J18").Value = 0.2521205402166373
$ws.Range("K18").Value = 3.192697466324603
$ws.Range("N18").Value = 2.561740972733048
$ws.Range("B19").Value = 0.1715551564166873
$ws.Range("D19").Value = 0.03144958750850435
$ws.Range("E19").Value = 0.1203772783564254
$ws.Range("F19").Value = 3.99846017212738
$ws.Range("G19").Value = 0.002575104017711711
$ws.Range("J19").Value = 0.2516212476845254
$ws.Range("K19").Value = 3.179684833473232
$ws.Range("N19").Value = 2.562775890957525
$ws.Range("B20").Value = 0.1755895349795225
$ws.Range("D20").Value = 0.0316757081531307
$ws.Range("E20").Value = 0.1214077851885378
$ws.Range("F20").Value = 4.018601594651386
$ws.Range("G20").Value = 0.002572975721641269
$ws.Range("J20").Value = 0.2538782003276907
$ws.Range("K20").Value = 3.238373595185408
$ws.Range("N20").Value = 2.558169704224696
$ws.Range("B21").Value = 0.1892314303709668
$ws.Range("D21").Value = 0.03245143305925779
$ws.Range("E21").Value = 0.1249500244191495
$ws.Range("F21").Value = 4.088853687904475
$ws.Range("G21").Value = 0.002566046112448055
$ws.Range("J21").Value = 0.2616385824592129
$ws.Range("K21").Value = 3.437859298199896
$ws.Range("N21").Value = 2.543589616582821
$ws.Range("B22").Value = 0.198211558276796
$ws.Range("D22").Value = 0.03297024127854797
$ws.Range("E22").Value = 0.1273249329201604
$ws.Range("F22").Value = 4.136703914485082
$ws.Range("G22").Value = 0.002561681191833802
$ws.Range("J22").Value = 0.2668437328067341
$ws.Range("K22").Value = 3.569968650073235
$ws.Range("N22").Value = 2.534731603907872
$ws.Range("B23").Value = 0.1934128336955752
$ws.Range("D23").Value = 0.03269226714149198
$ws.Range("E23").Value = 0.1260519025892393
$ws.Range("F23").Value = 4.110987272640443
$ws.Range("G23").Value = 0.002563995883553695
$ws.Range("J23").Value = 0.2640533768431936
$ws.Range("K23").Value = 3.499299446022633
$ws.Range("N23").Value = 2.539397704478645
$ws.Range("B24").Value = 0.1753676580159862
$ws.Range("D24").Value = 0.03166323029574869
$ws.Range("E24").Value = 0.1213508949695452
$ws.Range("F24").Value = 4.017485853841663
$ws.Range("G24").Value = 0.002573091763182312
$ws.Range("J24").Value = 0.2537535949854828
$ws.Range("K24").Value = 3.235142088499572
$ws.Range("N24").Value = 2.558419290700201
$ws.Range("B25").Value = 0.1561820874986353
$ws.Range("D25").Value = 0.03060461503518042
$ws.Range("E25").Value = 0.1165350047765124
$ws.Range("F25").Value = 3.924854191599252
$ws.Range("G25").Value = 0.002583612802405131
$ws.Range("J25").Value = 0.2432085963171176
$ws.Range("K25").Value = 2.95751981695139
$ws.Range("N25").Value = 2.581797932971853
